$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Riot Platforms, Inc. / RIOT) - minor value updates
$ws.Range("D2").Value = 14.84
$ws.Range("F2").Value = -8.029999999999999
$ws.Range("N2").Value = 50.68470204858703

# Row 3 (MARA Holdings, Inc. / MARA) - minor value updates
$ws.Range("E3").Value = 49.2
$ws.Range("F3").Value = 0.64
$ws.Range("N3").Value = 50.68470204858703

# Row 4 switches identity from Coinbase Global, Inc. / COIN to Bitcoin USD / BTC-USD
$ws.Range("B4").Value = "Bitcoin USD"
$ws.Range("C4").Value = "BTC-USD"
$ws.Range("D4").Value = 89214.23
$ws.Range("E4").Value = 58.1
$ws.Range("F4").Value = -1.31
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 63
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 48.4
$ws.Range("N4").Value = 50.68470204858703

# Row 5 switches identity from Bitcoin USD / BTC-USD to Coinbase Global, Inc. / COIN
$ws.Range("B5").Value = "Coinbase Global, Inc."
$ws.Range("C5").Value = "COIN"
$ws.Range("D5").Value = 268.25
$ws.Range("E5").Value = 43.6
$ws.Range("F5").Value = -1.67
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 53
$ws.Range("K5").Value = 48.2
$ws.Range("N5").Value = 50.68470204858703

# Row 6 (Strategy Inc / MSTR) - minor value updates
$ws.Range("D6").Value = 178.56
$ws.Range("E6").Value = 39.9
$ws.Range("F6").Value = 0.78
$ws.Range("G6").Value = 30
$ws.Range("K6").Value = 40.2
$ws.Range("N6").Value = 50.68470204858703
